$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59 (shifts existing rows 59-170 down to 60-171)
$ws.Rows("59:59").Insert()

# Populate the newly inserted row 59 with the new data entry
$ws.Range("A59").Value2 = 4
$ws.Range("B59").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C59").Value2 = "Los Lagos"
$ws.Range("D59").Value2 = 44533
$ws.Range("E59").Value2 = 10
$ws.Range("F59").Value2 = "Fruta"
$ws.Range("G59").Value2 = 100108
$ws.Range("H59").Value2 = "Tropicales y subtropicales"
$ws.Range("I59").Value2 = 100108005
$ws.Range("J59").Value2 = "Piña"
$ws.Range("K59").Value2 = "Caramelo"
$ws.Range("L59").Value2 = "Tercera"
$ws.Range("M59").Value2 = 200
$ws.Range("N59").Value2 = 20000
$ws.Range("O59").Value2 = 21000
$ws.Range("P59").Value2 = 20500
$ws.Range("Q59").Value2 = "$/caja 16 unidades"
$ws.Range("R59").Value2 = "Ecuador"
$ws.Range("S59").Value2 = 1281
$ws.Range("T59").Value2 = 16
